$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @("D2", "46.118.63"),
    @("E2", "  -1.38%  "),
    @("D3", "2.630.30"),
    @("E3", "  -0.05%  "),
    @("E4", "  -0.13%  "),
    @("D5", "310.48"),
    @("E5", "  -1.27%  "),
    @("D6", "98.90"),
    @("E6", "  -4.31%  "),
    @("D7", "0.597"),
    @("E7", "  -1.41%  "),
    @("E8", "  +0.03%  "),
    @("D9", "0.581"),
    @("E9", "  -2.05%  "),
    @("D10", "38.86"),
    @("E10", "  -1.29%  "),
    @("B11", "Dogecoin"),
    @("C11", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"),
    @("D11", "0.0846"),
    @("E11", "  -0.26%  "),
    @("B12", "OKB"),
    @("C12", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"),
    @("D12", "54.32"),
    @("E12", "  -1.33%  "),
    @("D13", "8.09"),
    @("E13", "  -3.00%  "),
    @("D14", "3.020.65"),
    @("E14", "  -0.40%  "),
    @("E15", "  +0.68%  "),
    @("D16", "2.623.73"),
    @("E16", "  -1.06%  "),
    @("D17", "0.921"),
    @("E17", "  -0.57%  "),
    @("D18", "14.91"),
    @("E18", "  -2.01%  "),
    @("D19", "46.103.61"),
    @("E19", "  -2.61%  "),
    @("E20", "  -1.04%  "),
    @("D21", "6.78"),
    @("E21", "  -0.56%  "),
    @("D22", "12.81"),
    @("E22", "  -4.05%  "),
    @("D23", "74.73"),
    @("E23", "  +4.66%  "),
    @("D24", "284.16"),
    @("E24", "  +9.55%  "),
    @("D25", "3.04"),
    @("E25", "  -2.69%  "),
    @("E26", "  +0.51%  "),
    @("D27", "30.21"),
    @("E27", "  +1.32%  "),
    @("D28", "1.00"),
    @("E28", "  +0.29%  "),
    @("B29", "Cosmos"),
    @("C29", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"),
    @("D29", "10.53"),
    @("E29", "  -1.82%  "),
    @("B30", "InjectiveProtocol"),
    @("C30", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"),
    @("D30", "38.79"),
    @("E30", "  -6.85%  "),
    @("B31", "Toncoin"),
    @("C31", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"),
    @("D31", "2.20"),
    @("E31", "  -4.30%  "),
    @("B32", "Filecoin"),
    @("C32", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"),
    @("D32", "6.25"),
    @("E32", "  -0.48%  "),
    @("B33", "ARBITRUM"),
    @("C33", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"),
    @("D33", "2.33"),
    @("E33", "  +1.02%  "),
    @("D34", "3.64"),
    @("E34", "  -3.32%  "),
    @("B35", "Monero"),
    @("C35", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"),
    @("D35", "156.73"),
    @("E35", "  +2.16%  "),
    @("B36", "Hedera"),
    @("C36", "https://coinranking.com/coin/jad286TjB+hedera-hbar"),
    @("D36", "0.0842"),
    @("E36", "  -0.46%  "),
    @("B37", "WEMIXToken"),
    @("C37", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"),
    @("D37", "2.82"),
    @("E37", "  -1.73%  "),
    @("B38", "Kaspa"),
    @("C38", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"),
    @("D38", "0.123"),
    @("E38", "  +3.84%  "),
    @("B39", "Stellar"),
    @("C39", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"),
    @("D39", "0.124"),
    @("E39", "  +0.52%  "),
    @("B40", "EnergySwap"),
    @("C40", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"),
    @("D40", "22.41"),
    @("E40", "  +3.64%  "),
    @("D41", "15.83"),
    @("E41", "  -7.17%  "),
    @("B42", "VeChain"),
    @("C42", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"),
    @("D42", "0.0329"),
    @("E42", "  -1.01%  "),
    @("B43", "NEARProtocol"),
    @("C43", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"),
    @("D43", "3.57"),
    @("E43", "  -4.26%  "),
    @("B44", "RenderToken"),
    @("C44", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"),
    @("D44", "4.04"),
    @("E44", "  -6.62%  "),
    @("B45", "Maker"),
    @("C45", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"),
    @("D45", "2.108.06"),
    @("E45", "  +3.30%  "),
    @("B46", "FirstDigitalUSD"),
    @("C46", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"),
    @("D46", "0.998"),
    @("E46", "  -0.04%  "),
    @("B47", "BitcoinSV"),
    @("C47", "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"),
    @("D47", "94.71"),
    @("E47", "  +1.77%  "),
    @("B48", "Aave"),
    @("C48", "https://coinranking.com/coin/ixgUfzmLR+aave-aave"),
    @("D48", "110.48"),
    @("E48", "  -3.31%  "),
    @("B49", "FraxShare"),
    @("C49", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"),
    @("D49", "9.15"),
    @("E49", "  -1.91%  "),
    @("B50", "RocketPoolETH"),
    @("C50", "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"),
    @("D50", "2.877.22"),
    @("E50", "  -0.61%  "),
    @("B51", "Algorand"),
    @("C51", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"),
    @("D51", "0.201"),
    @("E51", "  -1.36%  ")
)

foreach ($pair in $changes) {
    $addr = $pair[0]
    $val = $pair[1]
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}
